# Append the new daily COVID bulletin rows (2021-04-27 .. 2021-04-30)
# to "Planilha1" (sheet1), matching the commit "add data until April 30th".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Each entry: Date(serial), Descartados, EmInvestigacao, Confirmados, Examinados,
#             Recuperados, Ativos, Hospital, Domicilio, Obitos
$newRows = @(
    @(44313, 9201, 167, 3442, 12810, 3000, 361, 22, 339, 81),
    @(44314, 9249, 154, 3464, 12867, 3042, 340, 22, 318, 82),
    @(44315, 9317, 159, 3513, 12989, 3075, 354, 22, 332, 84),
    @(44316, 9349, 152, 3539, 13040, 3101, 353, 22, 331, 85)
)

$startRow = 300
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Move the frozen-pane selection to the new first empty row, as in the diff.
$lastRow = $startRow + $newRows.Count - 1
$nextRow = $lastRow + 1
[void]$ws.Range("A" + $nextRow).Select()
